$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'capri mens joggers'
$ws.Range("A2").Value = 'easton baseball pants mens'
$ws.Range("A3").Value = 'mueller knee pads basketball'
$ws.Range("A4").Value = 'bmx knee pads youth'
$ws.Range("A5").Value = 'basketball tights kids'
$ws.Range("A6").Value = 'tortoise knee pads'
$ws.Range("A7").Value = 'bendable knee pads'
$ws.Range("A8").Value = 'dakine knee pads'
$ws.Range("A9").Value = '187 knee pads'
$ws.Range("A10").Value = 'mma knee pads'
$ws.Range("A11").Value = 'scrubs men pants'
$ws.Range("A12").Value = 'copper compression pants'
$ws.Range("A13").Value = 'elbow pads knee pads'
$ws.Range("A14").Value = 'knee pad for scooter'
$ws.Range("A15").Value = 'hamstring compression pants'
$ws.Range("A16").Value = 'compression tights men basketball'
$ws.Range("A17").Value = 'white leggings for men'
$ws.Range("A18").Value = 'red leggings men'
$ws.Range("A19").Value = 'purple leggings men'
$ws.Range("A20").Value = 'yellow leggings men'
$ws.Range("A21").Value = 'athletic leggings for women'
$ws.Range("A22").Value = 'athletic leggings kids'
$ws.Range("A23").Value = 'basketball leggings for women'
$ws.Range("A24").Value = 'basketball leggings kids'
$ws.Range("A25").Value = 'bdu pants with knee pads'
$ws.Range("A26").Value = 'nike leggings for men'
$ws.Range("A27").Value = 'compression pants basketball'
$ws.Range("A28").Value = 'compression pants for women'
$ws.Range("A29").Value = 'compression pants kids'
$ws.Range("A30").Value = 'compression pants knee'
$ws.Range("A31").Value = 'compression pants men under armour'
$ws.Range("A32").Value = 'compression pants set'
$ws.Range("A33").Value = 'compression pants tesla'
$ws.Range("A34").Value = 'camo pants with knee pads'
$ws.Range("A35").Value = 'kids basketball pads'
$ws.Range("A36").Value = 'gold leggings men'
$ws.Range("A37").Value = 'kickboxing knee pads'
$ws.Range("A38").Value = 'men gym pants'
$ws.Range("A39").Value = 'men nike compression pants'
$ws.Range("A40").Value = 'swim leggings for men'
$ws.Range("A41").Value = 'fleece leggings men'
$ws.Range("A42").Value = 'sliding shorts with knee pads'
$ws.Range("A43").Value = 'mens basketball joggers'
$ws.Range("A44").Value = 'mcdavid basketball knee'
$ws.Range("A45").Value = '3xl knee pads'
$ws.Range("A46").Value = '511 knee pads'
$ws.Range("A47").Value = 'caving knee pads'
$ws.Range("A48").Value = 'enduro knee pads'
$ws.Range("A49").Value = 'bodyprox knee pads'
$ws.Range("A50").Value = 'blackhawk knee pads'
$ws.Range("A51").Value = 'snickers knee pads'
$ws.Range("A52").Value = 'bball knee pads'
$ws.Range("A53").Value = '661 knee pads'
$ws.Range("A54").Value = 'basket knee pads'
$ws.Range("A55").Value = 'armadillo knee pads'
$ws.Range("A56").Value = 'swim pants men'
$ws.Range("A57").Value = 'athletic capris'
$ws.Range("A58").Value = 'bunheads knee pads'
$ws.Range("A59").Value = '6xl compression pants'
$ws.Range("A60").Value = 'kali knee pads'
$ws.Range("A61").Value = 'arcteryx knee pads'
$ws.Range("A62").Value = 'spelunking knee pads'
$ws.Range("A63").Value = 'bcg compression pants'
$ws.Range("A64").Value = 'bcg knee pads'
$ws.Range("A65").Value = 'tesla mens leggings'
$ws.Range("A66").Value = 'awp knee pads'
$ws.Range("A67").Value = 'leatt knee pads'
$ws.Range("A68").Value = 'fr knee pads'
$ws.Range("A69").Value = 'eurotard knee pads'
$ws.Range("A70").Value = 'asics leggings men'
$ws.Range("A71").Value = 'alleson baseball pants youth'
$ws.Range("A72").Value = 'frozen knee pads'
$ws.Range("A73").Value = 'training tights men'
$ws.Range("A74").Value = 'training leggings for men'
$ws.Range("A75").Value = 'nba knee pads'
$ws.Range("A76").Value = 'prayer knee pads'
$ws.Range("A77").Value = 'neoprene knee pad'
$ws.Range("A78").Value = 's1 knee pads'
$ws.Range("A79").Value = 'white basketball pants'
$ws.Range("A80").Value = 'lululemon compression pants'
$ws.Range("A81").Value = 'xlarge knee pads'
$ws.Range("A82").Value = 'ballet tights for men'
$ws.Range("A83").Value = 'nike compression pants for men'
$ws.Range("A84").Value = 'mens compression pants adidas'
$ws.Range("A85").Value = 'mens compression pants champion'
$ws.Range("A86").Value = 'baseball pants mens with piping'
$ws.Range("A87").Value = 'exercise compression leggings'
$ws.Range("A88").Value = 'teen knee pads'
$ws.Range("A89").Value = 'short tights for men'
$ws.Range("A90").Value = 'runner compression pants'
$ws.Range("A91").Value = 'swimming tights men'
$ws.Range("A92").Value = 'half tights men'
$ws.Range("A93").Value = 'fishing knee pads'
$ws.Range("A94").Value = 'marvel knee pads'
$ws.Range("A95").Value = 'bauer compression pants'
$ws.Range("A96").Value = 'police knee pads'
$ws.Range("A97").Value = 'spiderman tights men'
$ws.Range("A98").Value = 'wwe knee pads'
$ws.Range("A99").Value = 'petite compression leggings'
$ws.Range("A100").Value = '3x compression leggings'
